$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: nickname, email, dienthoai
$ws.Range("E2").Value = "nam_abc_xyz"
$ws.Range("F2").Value = "nam_abc_xyz@gmail.com"
$ws.Range("G2").Value = "tỷyryrty"

# Remove the last two data rows (row 7 and row 8)
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()
